$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new (older) observation is inserted ahead of the existing series, so every
# existing data row (2-18) shifts down by one, to (3-19). Shift by copying
# the bounded A:E range row-by-row (bottom-up, so a source row is never
# clobbered before it's read) instead of a full-row Insert, which would drag
# the header row's bold formatting onto the vacated row and mint an unused
# style entry.
for ($r = 18; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":E" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":E" + ($r + 1))
    $src.Copy($dst)
}

# Populate the new first data row (oldest forecast origin/vectors) - reuses
# whatever formatting already lived in row 2 (same as every other data row).
$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 11.13090654781819
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 10.67037004222142

# The underlying simulation/evaluation was bugfixed, so every row's forecast
# columns (C = y_0_forecast, E = y_1_forecast) get refreshed values.
$forecasts = @(
    @(2,  11.13090654781819,   10.67037004222142),
    @(3,  4.672550446571067,   -0.7156496512470745),
    @(4,  -14.45332333832743,  7.857938327064184),
    @(5,  8.600536527919633,   12.64892828543749),
    @(6,  10.25770250047622,   10.40099841437159),
    @(7,  4.639893381363169,   6.662398279632087),
    @(8,  0.3058963467304165,  1.195213983078647),
    @(9,  4.068173739091874,   7.055025120039615),
    @(10, 4.984288257750213,   1.985659800779893),
    @(11, 1.878184267712912,   -0.3562142672005275),
    @(12, 4.695933104194339,   6.493919935864634),
    @(13, 4.892602738886098,   -2.576675125869599),
    @(14, 0.8049382522247184,  3.1919852842623),
    @(15, -8.784173899737169,  6.942816049735523),
    @(16, 5.110501195359984,   0.8094958705429534),
    @(17, 5.120680133083599,   0.5542886326586061),
    @(18, -0.5532735011319234, -3.561435976944571),
    @(19, -1.069674659641462,  0.01743232028155184)
)

foreach ($entry in $forecasts) {
    $r = $entry[0]
    $ws.Cells.Item($r, 3).Value = $entry[1]
    $ws.Cells.Item($r, 5).Value = $entry[2]
}
